$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 46, shifting existing rows 46:141 down to 47:142.
$ws.Rows.Item(46).Insert()

# Populate the newly inserted row 46 with the new record's data.
$ws.Range("A46").Value = 11
$ws.Range("B46").Value = "Vega Monumental Concepción"
$ws.Range("C46").Value = "Bíobío"
$ws.Range("D46").Value = 45203
$ws.Range("D46").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E46").Value = 8
$ws.Range("F46").Value = 100112012
$ws.Range("G46").Value = "Espinaca"
$ws.Range("H46").Value = "Sin especificar"
$ws.Range("I46").Value = "Primera"
$ws.Range("J46").Value = 50
$ws.Range("K46").Value = 11000
$ws.Range("L46").Value = 12000
$ws.Range("M46").Value = 11400
$ws.Range("N46").Value = "$/cuna 10 kilos"
$ws.Range("O46").Value = "Región Metropolitana"
$ws.Range("P46").Value = 1140
$ws.Range("Q46").Value = 10
$ws.Range("R46").Value = "Hortaliza"
